# Auto-generated: updates the cryptos price/volume table to the
# latest scrape (GitHub Actions refresh), including the two pairs of
# rows (40/41 and 44/45) whose coin name + link + values were swapped.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "39.707.15"
$ws.Range("E2").Value = "  -0.96%  "
$ws.Range("D3").Value = "2.190.47"
$ws.Range("E3").Value = "  -2.16%  "
$s = $ws.Range("D4").Style
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = $s
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "290.15"
$ws.Range("E5").Value = "  -1.22%  "
$ws.Range("D6").Value = "85.79"
$ws.Range("E6").Value = "  -1.59%  "
$ws.Range("D7").Value = "0.506"
$ws.Range("E7").Value = "  -1.96%  "
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("E9").Value = "  -2.90%  "
$ws.Range("D10").Value = "30.08"
$ws.Range("E10").Value = "  -3.71%  "
$ws.Range("D11").Value = "49.98"
$ws.Range("E11").Value = "  +6.41%  "
$ws.Range("E12").Value = "  -2.02%  "
$ws.Range("E13").Value = "  +2.35%  "
$ws.Range("E14").Value = "  -0.38%  "
$ws.Range("D15").Value = "2.529.30"
$ws.Range("E15").Value = "  -2.14%  "
$ws.Range("D16").Value = "13.64"
$ws.Range("E16").Value = "  -3.33%  "
$ws.Range("D17").Value = "2.148.22"
$ws.Range("E17").Value = "  -4.10%  "
$ws.Range("D18").Value = "0.723"
$ws.Range("E18").Value = "  -1.34%  "
$ws.Range("D19").Value = "39.597.18"
$ws.Range("E19").Value = "  -1.09%  "
$ws.Range("D20").Value = "0.0₃0879"
$ws.Range("E20").Value = "  -1.33%  "
$s = $ws.Range("D21").Style
$ws.Range("D21").Value = "'11.10"
$ws.Range("D21").Style = $s
$ws.Range("E21").Value = "  -1.15%  "
$ws.Range("E22").Value = "  -2.77%  "
$ws.Range("D23").Value = "64.94"
$ws.Range("E23").Value = "  -1.37%  "
$ws.Range("D24").Value = "235.98"
$ws.Range("E24").Value = "  -0.12%  "
$ws.Range("E25").Value = "  +0.20%  "
$ws.Range("E26").Value = "  -2.04%  "
$s = $ws.Range("D27").Style
$ws.Range("D27").Value = "'1.80"
$ws.Range("D27").Style = $s
$ws.Range("E27").Value = "  -3.07%  "
$ws.Range("D28").Value = "23.04"
$ws.Range("E28").Value = "  +0.31%  "
$ws.Range("E29").Value = "  -3.75%  "
$ws.Range("D30").Value = "9.13"
$ws.Range("E30").Value = "  -2.25%  "
$ws.Range("D31").Value = "155.56"
$ws.Range("E31").Value = "  +2.72%  "
$ws.Range("D32").Value = "31.09"
$ws.Range("E32").Value = "  -6.88%  "
$ws.Range("E33").Value = "  -0.05%  "
$s = $ws.Range("D34").Style
$ws.Range("D34").Value = "'4.90"
$ws.Range("D34").Style = $s
$ws.Range("E34").Value = "  -1.06%  "
$ws.Range("D35").Value = "0.0703"
$ws.Range("E35").Value = "  -2.66%  "
$ws.Range("E36").Value = "  -2.39%  "
$ws.Range("E37").Value = "  -0.26%  "
$ws.Range("D38").Value = "0.111"
$ws.Range("E38").Value = "  -0.26%  "
$ws.Range("D39").Value = "0.0969"
$ws.Range("E39").Value = "  -3.23%  "
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").Value = "1.67"
$ws.Range("E40").Value = "  -3.00%  "
$ws.Range("B41").Value = "Celestia"
$ws.Range("C41").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D41").Value = "15.04"
$ws.Range("E41").Value = "  -7.78%  "
$ws.Range("D42").Value = "2.112.41"
$ws.Range("E42").Value = "  +2.33%  "
$ws.Range("E43").Value = "  -3.67%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").Value = "0.0266"
$ws.Range("E44").Value = "  -1.28%  "
$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").Value = "2.08"
$ws.Range("E45").Value = "  -1.33%  "
$ws.Range("D46").Value = "9.67"
$ws.Range("E46").Value = "  -1.73%  "
$ws.Range("D47").Value = "17.17"
$ws.Range("E47").Value = "  -5.67%  "
$ws.Range("D48").Value = "2.65"
$ws.Range("E48").Value = "  +1.87%  "
$ws.Range("D49").Value = "2.394.27"
$ws.Range("E49").Value = "  -1.97%  "
$ws.Range("D50").Value = "1.46"
$ws.Range("E50").Value = "  +0.08%  "
$ws.Range("D51").Value = "87.75"
$ws.Range("E51").Value = "  -1.97%  "
